$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 12520.4
$ws.Range("J69").Value = 11324.75
$ws.Range("L69").Value = 33974.25
$ws.Range("N69").Value = -35722.25
$ws.Range("H72").Value = 12520.4
$ws.Range("J72").Value = 11324.75
$ws.Range("L72").Value = 101922.75
$ws.Range("N72").Value = -110658.75
$ws.Range("H80").Value = 1125
$ws.Range("I80").Value = 709.5
$ws.Range("J80").Value = 1309.6666
$ws.Range("K80").Value = 2128.5
$ws.Range("L80").Value = 3928.9998
$ws.Range("M80").Value = -1130.5
$ws.Range("N80").Value = -5924.9998
$ws.Range("H83").Value = 1125
$ws.Range("I83").Value = 709.5
$ws.Range("J83").Value = 1309.6666
$ws.Range("K83").Value = 6385.5
$ws.Range("L83").Value = 11786.9994
$ws.Range("M83").Value = -1393.5
$ws.Range("N83").Value = -21770.9994
$ws.Range("H103").Value = 844.89655
$ws.Range("J103").Value = 1717
$ws.Range("L103").Value = 5151
$ws.Range("N103").Value = -6323
$ws.Range("H106").Value = 5964.846
$ws.Range("I106").Value = 6303.1
$ws.Range("K106").Value = 6303.1
$ws.Range("M106").Value = -5672.1
$ws.Range("H115").Value = 279.3846
$ws.Range("I115").Value = 279.3846
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 838.1537999999999
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 728.8462000000001
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 7016.796
$ws.Range("I116").Value = 6870.2905
$ws.Range("J116").Value = 7269.1113
$ws.Range("K116").Value = 6870.2905
$ws.Range("L116").Value = 7269.1113
$ws.Range("M116").Value = -3428.2905
$ws.Range("N116").Value = -14153.1113
$ws.Range("H125").Value = 3710.0908
$ws.Range("I125").Value = 1929.5
$ws.Range("J125").Value = 5846.8
$ws.Range("K125").Value = 17365.5
$ws.Range("L125").Value = 52621.2
$ws.Range("M125").Value = -14905.5
$ws.Range("N125").Value = -57541.2
$ws.Range("H129").Value = 3033.5
$ws.Range("I129").Value = 2087.3
$ws.Range("K129").Value = 6261.900000000001
$ws.Range("M129").Value = -1261.900000000001
$ws.Range("H132").Value = 2569.2222
$ws.Range("I132").Value = 2537.8408
$ws.Range("K132").Value = 7613.5224
$ws.Range("M132").Value = -5083.5224
$ws.Range("H137").Value = 2192.3044
$ws.Range("I137").Value = 2029.3077
$ws.Range("J137").Value = 2404.2
$ws.Range("K137").Value = 6087.9231
$ws.Range("L137").Value = 7212.599999999999
$ws.Range("M137").Value = -3537.9231
$ws.Range("N137").Value = -12312.6
$ws.Range("H138").Value = 3144.606
$ws.Range("I138").Value = 3023.8
$ws.Range("K138").Value = 9071.400000000001
$ws.Range("M138").Value = -3931.400000000001
$ws.Range("H141").Value = 3152.9312
$ws.Range("I141").Value = 3034.12
$ws.Range("K141").Value = 9102.360000000001
$ws.Range("M141").Value = -3922.360000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 431.55
$ws.Range("I22").Value = 438.6875
$ws.Range("J22").Value = 403
$ws.Range("K22").Value = 438.6875
$ws.Range("L22").Value = 403
$ws.Range("M22").Value = -265.6875
$ws.Range("N22").Value = -749
$ws.Range("H32").Value = 38352.668
$ws.Range("J32").Value = 38352.668
$ws.Range("L32").Value = 38352.668
$ws.Range("N32").Value = -39120.668
$ws.Range("H107").Value = 8682.210999999999
$ws.Range("I107").Value = 9485.454
$ws.Range("K107").Value = 9485.454
$ws.Range("M107").Value = -7565.454

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9659
$ws.Range("I31").Value = 3012
$ws.Range("J31").Value = 10988.4
$ws.Range("K31").Value = 3012
$ws.Range("L31").Value = 10988.4
$ws.Range("M31").Value = -2717
$ws.Range("N31").Value = -11578.4
$ws.Range("H34").Value = 9659
$ws.Range("I34").Value = 3012
$ws.Range("J34").Value = 10988.4
$ws.Range("K34").Value = 3012
$ws.Range("L34").Value = 10988.4
$ws.Range("M34").Value = -2810
$ws.Range("N34").Value = -11392.4
$ws.Range("H58").Value = 3810.1667
$ws.Range("J58").Value = 2792.5
$ws.Range("L58").Value = 2792.5
$ws.Range("N58").Value = -3198.5
$ws.Range("H105").Value = 2159.25
$ws.Range("I105").Value = 2180.111
$ws.Range("K105").Value = 2180.111
$ws.Range("M105").Value = -433.1109999999999
$ws.Range("H122").Value = 4373.727
$ws.Range("I122").Value = 4373.727
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13121.181
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10671.181
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3477.3635
$ws.Range("I132").Value = 3376.3
$ws.Range("K132").Value = 10128.9
$ws.Range("M132").Value = -7598.900000000001
$ws.Range("H134").Value = 4131.5454
$ws.Range("I134").Value = 2770.5557
$ws.Range("K134").Value = 8311.667099999999
$ws.Range("M134").Value = -5776.667099999999
$ws.Range("H136").Value = 3810.1667
$ws.Range("J136").Value = 2792.5
$ws.Range("L136").Value = 8377.5
$ws.Range("N136").Value = -13477.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 9950
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 9950
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 29850
$ws.Range("N101").Value = -34718
$ws.Range("M101").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1488601.6
$ws.Range("I122").Value = 1783597.1
$ws.Range("J122").Value = 13624.25
$ws.Range("K122").Value = 5350791.300000001
$ws.Range("L122").Value = 40872.75
$ws.Range("M122").Value = -5348341.300000001
$ws.Range("N122").Value = -45772.75
$ws.Range("H126").Value = 3571.3
$ws.Range("I126").Value = 3571.3
$ws.Range("K126").Value = 10713.9
$ws.Range("M126").Value = -8243.900000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2941.1
$ws.Range("I7").Value = 2906.6924
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 2906.6924
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -2794.6924
$ws.Range("N7").Value = -3229
$ws.Range("H22").Value = 2512.7144
$ws.Range("I22").Value = 1099.6666
$ws.Range("J22").Value = 3572.5
$ws.Range("K22").Value = 1099.6666
$ws.Range("L22").Value = 3572.5
$ws.Range("M22").Value = -804.6666
$ws.Range("N22").Value = -4162.5
$ws.Range("H27").Value = 2512.7144
$ws.Range("I27").Value = 1099.6666
$ws.Range("J27").Value = 3572.5
$ws.Range("K27").Value = 1099.6666
$ws.Range("L27").Value = 3572.5
$ws.Range("M27").Value = -992.6666
$ws.Range("N27").Value = -3786.5
$ws.Range("H40").Value = 5572.4287
$ws.Range("I40").Value = 4928.8237
$ws.Range("K40").Value = 4928.8237
$ws.Range("M40").Value = -4792.8237
$ws.Range("H68").Value = 3312.2144
$ws.Range("J68").Value = 4829
$ws.Range("L68").Value = 4829
$ws.Range("N68").Value = -6327
$ws.Range("H71").Value = 3312.2144
$ws.Range("J71").Value = 4829
$ws.Range("L71").Value = 24145
$ws.Range("N71").Value = -31633
$ws.Range("H82").Value = 1766.8572
$ws.Range("I82").Value = 1810.4615
$ws.Range("K82").Value = 1810.4615
$ws.Range("M82").Value = -1449.4615
$ws.Range("H85").Value = 1766.8572
$ws.Range("I85").Value = 1810.4615
$ws.Range("K85").Value = 1810.4615
$ws.Range("M85").Value = -562.4614999999999
$ws.Range("H126").Value = 2941.1
$ws.Range("I126").Value = 2906.6924
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 8720.0772
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -6250.0772
$ws.Range("N126").Value = -13955
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 65181.09
$ws.Range("J76").Value = 67443.664
$ws.Range("L76").Value = 67443.664
$ws.Range("N76").Value = -68073.664
$ws.Range("H79").Value = 65181.09
$ws.Range("J79").Value = 67443.664
$ws.Range("L79").Value = 67443.664
$ws.Range("N79").Value = -69627.664
